# Auto-committed on 2023/09/15 週五 17:07:32.90
#
# The DBD sheet's field-type list used the literal "DATE" type for the two
# timestamp columns (CreateDate / LastUpdate). Replace it with "TIMESTAMP".

$wb = $excel.ActiveWorkbook

$wsDBD = $wb.Worksheets.Item("DBD")
$wsDBS = $wb.Worksheets.Item("DBS")

# D13 -> 建檔日期時間 / CreateDate row ; D15 -> 最後更新日期時間 / LastUpdate row
$wsDBD.Range("D13").Value = "TIMESTAMP"
$wsDBD.Range("D15").Value = "TIMESTAMP"

# Those two rows had a leftover manual row-height override (wrapped text from
# an earlier edit); re-autofit them back down to the sheet's default height.
$wsDBD.Rows.Item(14).AutoFit()
$wsDBD.Rows.Item(16).AutoFit()

# Restore DBS's own (no-longer-selected) cursor position first …
$wsDBS.Activate()
$wsDBS.Range("A7").Select()

# … then make DBD the active sheet/tab and set its own selection, matching
# the re-saved view state (DBD becomes the selected/visible tab).
$wsDBD.Activate()
$wsDBD.Range("D13").Select()
